$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<your>"
$ws.Range("C2").Value = 37

# Row 4
$ws.Range("C4").Value = 38

# Row 5
$ws.Range("B5").Value = "<yes>"

# Row 7
$ws.Range("B7").Value = "<otim>"
$ws.Range("C7").Value = 35

# Row 8
$ws.Range("C8").Value = 28

# Row 9
$ws.Range("C9").Value = 46

# Row 10
$ws.Range("B10").Value = "<on>"
$ws.Range("C10").Value = 33

# Row 11
$ws.Range("B11").Value = "<been>"
$ws.Range("C11").Value = 27

# Row 12
$ws.Range("C12").Value = 27

# Row 13
$ws.Range("C13").Value = 35

# Row 14
$ws.Range("B14").Value = "<they>"

# Row 15
$ws.Range("C15").Value = 28

# Row 16
$ws.Range("C16").Value = 36

# Row 17
$ws.Range("C17").Value = 31

# Row 18
$ws.Range("C18").Value = 33
